$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): unchanged text, refreshed to new shared-string run ---
$ws.Cells.Item(1,1).Value = 'Row'
$ws.Cells.Item(1,2).Value = 'Prognose'
$ws.Cells.Item(1,3).Value = 'surveys'
$ws.Cells.Item(1,4).Value = 'production'
$ws.Cells.Item(1,5).Value = 'orders'
$ws.Cells.Item(1,6).Value = 'turnover'
$ws.Cells.Item(1,7).Value = 'financial'
$ws.Cells.Item(1,8).Value = 'labor market'
$ws.Cells.Item(1,9).Value = 'prices'
$ws.Cells.Item(1,10).Value = 'national accounts'
$ws.Cells.Item(1,11).Value = 'Revision'

# --- Column A (dates) forced to text, rows 2-12 ---
$defaultStyle = $ws.Cells.Item(1,1).Style
$dateCol = @('2025-03-30', '2025-04-15', '2025-04-30', '2025-05-15', '2025-05-30', '2025-06-15', '2025-06-30', '2025-07-15', '2025-07-30', '2025-08-15', '2025-08-30')
for ($i = 0; $i -lt $dateCol.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dateCol[$i]
    $cell.Style = $defaultStyle
}

# --- Numeric data, columns B-K, rows 2-12 ---
$ws.Cells.Item(2,2).Value = 0.2911636164037934
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 0
$ws.Cells.Item(3,2).Value = 0.29136534467023101
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = [double]"-4.7098465475294568e-05"
$ws.Cells.Item(3,5).Value = [double]"9.9459169207487835e-05"
$ws.Cells.Item(3,6).Value = [double]"5.0235990409296673e-06"
$ws.Cells.Item(3,7).Value = [double]"8.6029524228423345e-06"
$ws.Cells.Item(3,8).Value = [double]"-1.8004373285960302e-06"
$ws.Cells.Item(3,9).Value = [double]"3.5577259129439337e-05"
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).Value = [double]"-7.4250960710287028e-05"
$ws.Cells.Item(4,2).Value = 0.28967452880841077
$ws.Cells.Item(4,3).Value = -0.00052809552062833527
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = [double]"1.1403316501113045e-06"
$ws.Cells.Item(4,6).Value = [double]"3.2425268034234206e-06"
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = [double]"4.288258152585392e-06"
$ws.Cells.Item(4,9).Value = [double]"1.6426382003275625e-05"
$ws.Cells.Item(4,10).Value = [double]"-3.9673963377687186e-05"
$ws.Cells.Item(4,11).Value = [double]"4.5647522230907178e-05"
$ws.Cells.Item(5,2).Value = 0.28691914465085133
$ws.Cells.Item(5,3).Value = 0.00041142587951390792
$ws.Cells.Item(5,4).Value = -0.0005355935051578717
$ws.Cells.Item(5,5).Value = -0.00015120801069230279
$ws.Cells.Item(5,6).Value = -0.00042953206308585516
$ws.Cells.Item(5,7).Value = -0.00017285228764353869
$ws.Cells.Item(5,8).Value = [double]"-1.5400701604121397e-05"
$ws.Cells.Item(5,9).Value = [double]"-3.5270070739567232e-05"
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = [double]"5.0266491258221802e-05"
$ws.Cells.Item(6,2).Value = 0.30199844494534611
$ws.Cells.Item(6,3).Value = 0.0068055803533970151
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = [double]"4.3460976503154901e-07"
$ws.Cells.Item(6,6).Value = [double]"-5.1116093642896521e-05"
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(6,8).Value = [double]"-1.5343767745636981e-05"
$ws.Cells.Item(6,9).Value = -0.00032455457911904862
$ws.Cells.Item(6,10).Value = 0
$ws.Cells.Item(6,11).Value = -0.00069136854246054646
$ws.Cells.Item(7,2).Value = 0.29673150501952328
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = -0.0014172882290721873
$ws.Cells.Item(7,5).Value = [double]"-6.5539608027155748e-05"
$ws.Cells.Item(7,6).Value = -0.00064854301684427306
$ws.Cells.Item(7,7).Value = [double]"7.0007204657696818e-05"
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = [double]"3.6662663963353803e-05"
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = -0.0005009820080352223
$ws.Cells.Item(8,2).Value = 0.3174168217283419
$ws.Cells.Item(8,3).Value = 0.0083044396698239743
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = [double]"-6.535629445892746e-07"
$ws.Cells.Item(8,6).Value = -0.00028123860817182199
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = [double]"7.5755858353043634e-07"
$ws.Cells.Item(8,9).Value = -0.00047470596013378001
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = [double]"1.6407094301729153e-05"
$ws.Cells.Item(9,2).Value = 0.30017669794544355
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = 0.00036872819708495564
$ws.Cells.Item(9,5).Value = -0.0011390029640799596
$ws.Cells.Item(9,6).Value = -0.0061574420888025425
$ws.Cells.Item(9,7).Value = [double]"5.1541026403928984e-05"
$ws.Cells.Item(9,8).Value = -0.00020156826556262826
$ws.Cells.Item(9,9).Value = [double]"-1.8768452536362621e-05"
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = [double]"4.8534347821038715e-05"
$ws.Cells.Item(10,2).Value = 0.2744272584644743
$ws.Cells.Item(10,3).Value = 0.0021685996602088709
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = [double]"7.6147396097620782e-07"
$ws.Cells.Item(10,6).Value = -0.0001938541292178688
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = [double]"3.467196603941974e-06"
$ws.Cells.Item(10,9).Value = [double]"-4.3495937542944358e-05"
$ws.Cells.Item(10,10).Value = -0.0022210865503780354
$ws.Cells.Item(10,11).Value = 0.00082417556161812344
$ws.Cells.Item(11,2).Value = 0.26537080936878038
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0.0031170447933382165
$ws.Cells.Item(11,5).Value = -0.00066696701731159159
$ws.Cells.Item(11,6).Value = -0.0041151407683819993
$ws.Cells.Item(11,7).Value = [double]"1.0298861605260475e-05"
$ws.Cells.Item(11,8).Value = [double]"-2.4682683877119419e-05"
$ws.Cells.Item(11,9).Value = [double]"-4.8274579780456181e-05"
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).Value = 0.0011938580748458438
$ws.Cells.Item(12,2).Value = 0.21998524273585271
$ws.Cells.Item(12,3).Value = -0.034424313886733074
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = [double]"7.4350782594187979e-06"
$ws.Cells.Item(12,6).Value = [double]"1.1677220630371584e-06"
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = [double]"-1.4008127900192719e-06"
$ws.Cells.Item(12,9).Value = 0.00038163294133230683
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,11).Value = -0.0037435821322276963
